$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($rng, $val)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue $ws.Range("D2") "304.18"
Set-TextValue $ws.Range("E2") "1.29%"
Set-TextValue $ws.Range("G2") "4"
Set-TextValue $ws.Range("E3") "-5.43%"
Set-TextValue $ws.Range("G3") "4"
Set-TextValue $ws.Range("D4") "5.022"
Set-TextValue $ws.Range("E4") "0.88%"
Set-TextValue $ws.Range("G4") "4"
Set-TextValue $ws.Range("D5") "0.07808"
Set-TextValue $ws.Range("E5") "1.03%"
Set-TextValue $ws.Range("G5") "4"
Set-TextValue $ws.Range("D6") "2.173"
Set-TextValue $ws.Range("E6") "-0.17%"
Set-TextValue $ws.Range("G6") "4"
Set-TextValue $ws.Range("D7") "7.904"
Set-TextValue $ws.Range("E7") "-0.73%"
Set-TextValue $ws.Range("G7") "4"
Set-TextValue $ws.Range("D8") "4.096"
Set-TextValue $ws.Range("E8") "2.38%"
Set-TextValue $ws.Range("G8") "4"
Set-TextValue $ws.Range("D9") "0.9172"
Set-TextValue $ws.Range("E9") "0.33%"
Set-TextValue $ws.Range("G9") "4"
Set-TextValue $ws.Range("D10") "0.09721"
Set-TextValue $ws.Range("E10") "7.14%"
Set-TextValue $ws.Range("G10") "4"
Set-TextValue $ws.Range("D11") "0.1859"
Set-TextValue $ws.Range("E11") "3.42%"
Set-TextValue $ws.Range("G11") "4"
Set-TextValue $ws.Range("D12") "0.08601"
Set-TextValue $ws.Range("E12") "2.29%"
Set-TextValue $ws.Range("G12") "4"
Set-TextValue $ws.Range("D13") "0.03496"
Set-TextValue $ws.Range("E13") "-0.87%"
Set-TextValue $ws.Range("G13") "4"
Set-TextValue $ws.Range("D14") "0.09913"
Set-TextValue $ws.Range("E14") "-0.28%"
Set-TextValue $ws.Range("G14") "4"
Set-TextValue $ws.Range("D15") "0.001430"
Set-TextValue $ws.Range("E15") "-3.64%"
Set-TextValue $ws.Range("G15") "4"
Set-TextValue $ws.Range("D16") "0.005673"
Set-TextValue $ws.Range("E16") "-0.35%"
Set-TextValue $ws.Range("G16") "4"
Set-TextValue $ws.Range("D17") "3.457"
Set-TextValue $ws.Range("E17") "-0.53%"
Set-TextValue $ws.Range("G17") "4"
Set-TextValue $ws.Range("D18") "2.393"
Set-TextValue $ws.Range("E18") "7.65%"
Set-TextValue $ws.Range("G18") "4"
Set-TextValue $ws.Range("D19") "0.3425"
Set-TextValue $ws.Range("E19") "-1.09%"
Set-TextValue $ws.Range("G19") "4"
Set-TextValue $ws.Range("E20") "2.22%"
Set-TextValue $ws.Range("G20") "4"
Set-TextValue $ws.Range("D21") "4.773"
Set-TextValue $ws.Range("E21") "4.91%"
Set-TextValue $ws.Range("G21") "4"
Set-TextValue $ws.Range("D22") "0.2210"
Set-TextValue $ws.Range("E22") "-1.02%"
Set-TextValue $ws.Range("G22") "4"
Set-TextValue $ws.Range("D23") "0.04598"
Set-TextValue $ws.Range("E23") "-1.39%"
Set-TextValue $ws.Range("G23") "4"
Set-TextValue $ws.Range("E24") "14.83%"
Set-TextValue $ws.Range("G24") "4"
Set-TextValue $ws.Range("D25") "0.001230"
Set-TextValue $ws.Range("E25") "0.08%"
Set-TextValue $ws.Range("G25") "4"
Set-TextValue $ws.Range("D26") "0.0001400"
Set-TextValue $ws.Range("E26") "7.58%"
Set-TextValue $ws.Range("G26") "4"
Set-TextValue $ws.Range("E27") "0.01%"
Set-TextValue $ws.Range("G27") "4"
Set-TextValue $ws.Range("G28") "4"
Set-TextValue $ws.Range("G29") "4"
Set-TextValue $ws.Range("G30") "4"
Set-TextValue $ws.Range("G31") "4"
Set-TextValue $ws.Range("G32") "4"
Set-TextValue $ws.Range("G33") "4"
Set-TextValue $ws.Range("G34") "4"
Set-TextValue $ws.Range("G35") "4"
Set-TextValue $ws.Range("G36") "4"
Set-TextValue $ws.Range("G37") "4"
Set-TextValue $ws.Range("G38") "4"
Set-TextValue $ws.Range("D39") "0.01824"
Set-TextValue $ws.Range("E39") "4.59%"
Set-TextValue $ws.Range("G39") "4"
Set-TextValue $ws.Range("D40") "0.04723"
Set-TextValue $ws.Range("E40") "0.97%"
Set-TextValue $ws.Range("G40") "4"
Set-TextValue $ws.Range("D41") "0.007472"
Set-TextValue $ws.Range("E41") "-7.52%"
Set-TextValue $ws.Range("G41") "4"
Set-TextValue $ws.Range("D42") "0.1396"
Set-TextValue $ws.Range("E42") "0.64%"
Set-TextValue $ws.Range("G42") "4"
Set-TextValue $ws.Range("D43") "0.007753"
Set-TextValue $ws.Range("E43") "1.04%"
Set-TextValue $ws.Range("G43") "4"
Set-TextValue $ws.Range("D44") "0.002231"
Set-TextValue $ws.Range("E44") "-3.14%"
Set-TextValue $ws.Range("G44") "4"
Set-TextValue $ws.Range("D45") "0.01109"
Set-TextValue $ws.Range("E45") "9.56%"
Set-TextValue $ws.Range("G45") "4"
Set-TextValue $ws.Range("D46") "0.00006359"
Set-TextValue $ws.Range("E46") "5.61%"
Set-TextValue $ws.Range("G46") "4"
Set-TextValue $ws.Range("E47") "-0.03%"
Set-TextValue $ws.Range("G47") "4"
Set-TextValue $ws.Range("E48") "0.07%"
Set-TextValue $ws.Range("G48") "4"
Set-TextValue $ws.Range("D49") "25.50"
Set-TextValue $ws.Range("E49") "195.07%"
Set-TextValue $ws.Range("G49") "4"
Set-TextValue $ws.Range("G50") "4"
Set-TextValue $ws.Range("E51") "-0.03%"
Set-TextValue $ws.Range("G51") "4"
